# scaling-ob.xlsx — add the "EDI" / "SCOT (Scaled)" row labels in column A
# for the two data blocks on Sheet1 (rows 2-9 and 12-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label cells (these introduce the two new shared strings "EDI" and
# "SCOT (Scaled)" and bump A2/A12 from empty to labelled rows).
$ws.Range("A2").Value  = "EDI"
$ws.Range("A12").Value = "SCOT (Scaled)"

# Column A now holds real text, so it needs to widen to fit ("SCOT (Scaled)"
# is the longest entry) instead of sharing the default sheet column width.
$ws.Columns.Item(1).ColumnWidth = 11.75

# Leave the selection where the author left it after making the edit.
$ws.Range("K21").Select() | Out-Null
